$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct capitalization of the existing "Product 3" description
$ws.Range("B3").Value = "The third and final one"

# Add the fourth sample product as a new row
$ws.Range("A4").Value = "Product 4"
$ws.Range("B4").Value = "Testing css overflow"
$ws.Range("C4").Value = 23
$ws.Range("D4").Value = "P4"

# Leave the selection on the last cell that was entered
$null = $ws.Range("D4").Select()
